$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The three header columns "擬調" (O), "下限" (P) and "上限" (Q) were removed
# from the worksheet; deleting the entire columns shifts the former column R
# ("調整後利率") left into column O.
$ws.Range("O1:Q1").EntireColumn.Delete()

# Keep the workbook's hidden _FilterDatabase defined name in sync with the
# worksheet's new used range (was $A$1:$R$1, now $A$1:$O$1).
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=正常件!`$A`$1:`$O`$1"
    }
}

# Restore the cursor position recorded after the edit.
[void]$ws.Range("M8").Select()
